# Add a new slide ("Disadvantages") at the end of the deck, using the
# same "Title and Content" layout (index 2) as the other slides.
$p = $ppt.ActivePresentation
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Disadvantages"

# Body / content placeholder - one paragraph per bullet point.
$bodyLines = @(
    "Slower training on large datasets",
    "Requires clean, fully preprocessed data",
    "Easily confused by incorrect or unusual data",
    "Struggles with imbalanced datasets",
    "Performs best with simple models and overfits with complex ones",
    "Less accurate compared to XGBoost and LGBoost algorithms"
)
$contentShape = $s.Shapes.Item(2)
$contentShape.TextFrame.TextRange.Text = [string]::Join("`r", $bodyLines)

# Reposition/resize the content placeholder to match the authored slide
# (values below are the point-equivalents of the target EMU offsets,
# chosen so the single-precision round trip lands back on the exact
# EMU integers: 1186056 / 2995386 / 8825659 / 3416300).
$contentShape.Left = 93.39024622047243
$contentShape.Top = 235.8571753543307
$contentShape.Width = 694.9338095275591
$contentShape.Height = 269.0
